$d = $word.ActiveDocument

$d.Content.Find.Execute("136÷2=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "397÷8=49, 5", 2) | Out-Null
$d.Content.Find.Execute("537÷5=107, 2", $true, $false, $false, $false, $false, $true, 1, $false, "331÷3=110, 1", 2) | Out-Null
$d.Content.Find.Execute("154÷7=22, 0", $true, $false, $false, $false, $false, $true, 1, $false, "926÷7=132, 2", 2) | Out-Null
$d.Content.Find.Execute("905÷3=301, 2", $true, $false, $false, $false, $false, $true, 1, $false, "472÷3=157, 1", 2) | Out-Null
$d.Content.Find.Execute("204÷2=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "173÷9=19, 2", 2) | Out-Null
$d.Content.Find.Execute("225÷8=28, 1", $true, $false, $false, $false, $false, $true, 1, $false, "532÷3=177, 1", 2) | Out-Null
$d.Content.Find.Execute("110÷7=15, 5", $true, $false, $false, $false, $false, $true, 1, $false, "563÷7=80, 3", 2) | Out-Null
$d.Content.Find.Execute("151÷4=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "683÷7=97, 4", 2) | Out-Null
$d.Content.Find.Execute("717÷9=79, 6", $true, $false, $false, $false, $false, $true, 1, $false, "287÷8=35, 7", 2) | Out-Null
$d.Content.Find.Execute("560÷5=112, 0", $true, $false, $false, $false, $false, $true, 1, $false, "167÷2=83, 1", 2) | Out-Null
$d.Content.Find.Execute("222÷9=24, 6", $true, $false, $false, $false, $false, $true, 1, $false, "918÷8=114, 6", 2) | Out-Null
$d.Content.Find.Execute("120÷8=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "243÷8=30, 3", 2) | Out-Null
$d.Content.Find.Execute("308÷8=38, 4", $true, $false, $false, $false, $false, $true, 1, $false, "196÷5=39, 1", 2) | Out-Null
$d.Content.Find.Execute("755÷9=83, 8", $true, $false, $false, $false, $false, $true, 1, $false, "853÷8=106, 5", 2) | Out-Null
$d.Content.Find.Execute("201÷6=33, 3", $true, $false, $false, $false, $false, $true, 1, $false, "511÷2=255, 1", 2) | Out-Null
$d.Content.Find.Execute("396÷3=132, 0", $true, $false, $false, $false, $false, $true, 1, $false, "888÷4=222, 0", 2) | Out-Null
$d.Content.Find.Execute("430÷9=47, 7", $true, $false, $false, $false, $false, $true, 1, $false, "325÷4=81, 1", 2) | Out-Null
$d.Content.Find.Execute("473÷3=157, 2", $true, $false, $false, $false, $false, $true, 1, $false, "475÷3=158, 1", 2) | Out-Null
$d.Content.Find.Execute("257÷3=85, 2", $true, $false, $false, $false, $false, $true, 1, $false, "979÷6=163, 1", 2) | Out-Null
$d.Content.Find.Execute("253÷4=63, 1", $true, $false, $false, $false, $false, $true, 1, $false, "212÷3=70, 2", 2) | Out-Null
$d.Content.Find.Execute("165÷5=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "369÷8=46, 1", 2) | Out-Null
$d.Content.Find.Execute("142÷7=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "414÷9=46, 0", 2) | Out-Null
$d.Content.Find.Execute("835÷3=278, 1", $true, $false, $false, $false, $false, $true, 1, $false, "173÷7=24, 5", 2) | Out-Null
$d.Content.Find.Execute("230÷8=28, 6", $true, $false, $false, $false, $false, $true, 1, $false, "881÷9=97, 8", 2) | Out-Null
$d.Content.Find.Execute("276÷7=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "691÷5=138, 1", 2) | Out-Null
